$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Refresh the "time_taken" timestamps (column F) on the "data" sheet
# ---------------------------------------------------------------------------
$ws1.Range("F2").Value = "2021-10-05 14:34:01.846217"
$ws1.Range("F3").Value = "2021-10-05 14:34:01.846226"
$ws1.Range("F4").Value = "2021-10-05 14:34:01.846230"
$ws1.Range("F5").Value = "2021-10-05 14:34:01.846233"
$ws1.Range("F6").Value = "2021-10-05 14:34:01.846236"
$ws1.Range("F7").Value = "2021-10-05 14:34:01.846239"
$ws1.Range("F8").Value = "2021-10-05 14:34:01.846242"
$ws1.Range("F9").Value = "2021-10-05 14:34:01.846245"
$ws1.Range("F10").Value = "2021-10-05 14:34:01.846248"
$ws1.Range("F11").Value = "2021-10-05 14:34:01.846251"
$ws1.Range("F12").Value = "2021-10-05 14:34:01.846254"
$ws1.Range("F13").Value = "2021-10-05 14:34:01.846256"
$ws1.Range("F14").Value = "2021-10-05 14:34:01.846260"
$ws1.Range("F15").Value = "2021-10-05 14:34:01.846262"
$ws1.Range("F16").Value = "2021-10-05 14:34:01.846265"
$ws1.Range("F17").Value = "2021-10-05 14:34:01.846268"
$ws1.Range("F18").Value = "2021-10-05 14:34:01.846271"
$ws1.Range("F19").Value = "2021-10-05 14:34:01.846274"
$ws1.Range("F20").Value = "2021-10-05 14:34:01.846277"
$ws1.Range("F21").Value = "2021-10-05 14:34:01.846280"
$ws1.Range("F22").Value = "2021-10-05 14:34:01.846283"
$ws1.Range("F23").Value = "2021-10-05 14:34:01.846286"
$ws1.Range("F24").Value = "2021-10-05 14:34:01.846288"
$ws1.Range("F25").Value = "2021-10-05 14:34:01.846291"
$ws1.Range("F26").Value = "2021-10-05 14:34:01.846294"
$ws1.Range("F27").Value = "2021-10-05 14:34:01.846297"
$ws1.Range("F28").Value = "2021-10-05 14:34:01.846300"

# ---------------------------------------------------------------------------
# 2. Add a new "metadata" worksheet directly after "data"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "metadata"

# Match the page margins used by the "data" sheet (0.75/0.75/1/1/0.5/0.5 in)
$ws2.PageSetup.LeftMargin = 54
$ws2.PageSetup.RightMargin = 54
$ws2.PageSetup.TopMargin = 72
$ws2.PageSetup.BottomMargin = 72
$ws2.PageSetup.HeaderMargin = 36
$ws2.PageSetup.FooterMargin = 36

# Re-use the header look (bold font + thin border) from the "data" sheet
$ws1.Range("B1:F1").Copy()
$ws2.Range("B1:F1").PasteSpecial(-4122)
$ws1.Range("B1").Copy()
$ws2.Range("G1").PasteSpecial(-4122)

# Re-use the index-column look from the "data" sheet for A2
$ws1.Range("A2").Copy()
$ws2.Range("A2").PasteSpecial(-4122)

# Header row
$ws2.Range("B1").Value = "data_name"
$ws2.Range("C1").Value = "data_id"
$ws2.Range("D1").Value = "data_version"
$ws2.Range("E1").Value = "data_version_created"
$ws2.Range("F1").Value = "panel_query_time"
$ws2.Range("G1").Value = "panel_get_request"

# Data row
$ws2.Range("A2").Value = 0
$ws2.Range("B2").Value = "Holoprosencephaly and septo-optic dysplasia"
$ws2.Range("C2").Value = 112

# "data_version" must stay textual ("1.2"), not become the number 1.2
$d2 = $ws2.Range("D2")
$d2.NumberFormat = "@"
$d2.Value = "1.2"
$d2.Style = "Normal"

$ws2.Range("E2").Value = "2021-09-15T00:39:48.042033Z"
$ws2.Range("F2").Value = "2021-10-05 14:34:01.842963"
$ws2.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/112/?format=json"
